$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2,3 down to 3,4)
$ws.Rows("2:2").Insert()

$A2 = '莎普爱思滴眼睛'

$B2 = @'
莎普爱思滴眼睛	用词模糊,消息虚假,夸大失实	白内障，看不清~莎普爱思滴眼睛~
白内障，看不清~莎普爱思滴滴滴！
模糊滴！重影滴！黑影滴！
白内障，看不清~莎普爱思滴眼睛~
白内障，看不清~莎普爱思滴滴滴！
模糊滴！重影滴！黑影滴！有点痛！坚持滴！！！。 
'@

$C2 = '用词模糊,消息虚假,夸大失实'
$D2 = '我帮你看了一下哈～这个广告使用了模糊的用词和夸大的手法，来宣传莎普爱思滴眼液对白内障的治疗效果。虽然滴眼液可能对某些人有帮助，但这种夸大和不清晰的表述可能存在误导性。建议你谨慎对待这样的广告宣传，并且在购买前咨询专业的医生或药剂师，获取更准确的产品信息和使用建议。避免购买🚫'
$E2 = '检索到减肥茶和糖尿病的知识，输出完全错误'
$F2 = '知识库'

$ws.Range("A2").Value = $A2
$ws.Range("B2").Value = $B2
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = $D2
$ws.Range("E2").Value = $E2
$ws.Range("F2").Value = $F2
